# Add files via upload
# Append 17 new Spanish/Kawaiinese vocabulary rows (152-168) to the single
# worksheet, growing the used range from A1:C151 to A1:C168, and move the
# sheet's view/selection down to the newly added data, mirroring the
# upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A = Spanish, Column B = Kawaiinese.
# Assignment order below reproduces the exact order new vocabulary entries
# were typed in (e.g. rows 159/160 had their Spanish column filled before
# either Kawaiinese cell, and row 166 had column B typed before column A).
$ws.Cells.Item(152,1).Value = "Lugar"
$ws.Cells.Item(152,2).Value = "Pureisu"

$ws.Cells.Item(153,1).Value = "Cuyo / El cual"
$ws.Cells.Item(153,2).Value = "Utu"

$ws.Cells.Item(154,1).Value = "Preciso"
$ws.Cells.Item(154,2).Value = "Purusisu"

$ws.Cells.Item(155,1).Value = "Nombre"
$ws.Cells.Item(155,2).Value = "Denomu"

$ws.Cells.Item(156,1).Value = "Hace poco tiempo"
$ws.Cells.Item(156,2).Value = "Nouso"

$ws.Cells.Item(157,1).Value = "Un"
$ws.Cells.Item(157,2).Value = "a"

$ws.Cells.Item(158,1).Value = "Hidalgo"
$ws.Cells.Item(158,2).Value = "Nobaha"

$ws.Cells.Item(159,1).Value = "Norumi"
$ws.Cells.Item(160,1).Value = "Norumimeru"
$ws.Cells.Item(159,2).Value = "Normal"
$ws.Cells.Item(160,2).Value = "Normalmente"

$ws.Cells.Item(161,1).Value = "Llevar"
$ws.Cells.Item(161,2).Value = "Kipu"

$ws.Cells.Item(162,1).Value = "Lanza"
$ws.Cells.Item(162,2).Value = "Paropu"

$ws.Cells.Item(163,1).Value = "Escudo"
$ws.Cells.Item(163,2).Value = "Esukudo"

$ws.Cells.Item(164,1).Value = "Viejo"
$ws.Cells.Item(164,2).Value = "Orudo"

$ws.Cells.Item(165,1).Value = "Antiguo"
$ws.Cells.Item(165,2).Value = "Orudo"

$ws.Cells.Item(166,2).Value = "Bariia"
$ws.Cells.Item(166,1).Value = "Flaco"

$ws.Cells.Item(167,1).Value = "Caballo"
$ws.Cells.Item(167,2).Value = "Hibaro"

$ws.Cells.Item(168,1).Value = "Perro"
$ws.Cells.Item(168,2).Value = "Dogu"

# Scroll the view down to the new rows and select the last edited cell,
# matching the saved view state of the uploaded workbook.
[void]$ws.Range("B168").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 148
$win.ScrollColumn = 1
